$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing "204..215" numbering (column A) / near-zero residual
# series (column B) by 12 more rows (rows 206:217), matching the same
# formatting (bold/border/center style carried by A2:A205) as the rest of
# the column-A index.
$ws.Range("A205").Copy()
$ws.Range("A206:A217").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$aValues = @(204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215)
$bValues = @(
    -[double]"8.326672684688674E-17",
    0,
    [double]"4.996003610813205E-17",
    [double]"1.657624654822282E-17",
    [double]"1.734723475976807E-18",
    [double]"6.046750401976298E-17",
    -[double]"5.204170427930421E-18",
    [double]"2.914335439641036E-17",
    0,
    -[double]"5.088522196198634E-17",
    0,
    0
)

$row = 206
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $row = $row + 1
}
